$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 213, shifting existing rows 213..269 down to 214..270
$ws.Rows.Item(213).EntireRow.Insert()

# Populate the newly inserted row 213 with the new data record
$ws.Range("A213").Value = 6
$ws.Range("B213").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C213").Value = "Metropolitana"
$ws.Range("D213").Value = 44855
$ws.Range("E213").Value = 13
$ws.Range("F213").Value = 100112022
$ws.Range("G213").Value = "Arveja Verde"
$ws.Range("H213").Value = "Sin especificar"
$ws.Range("I213").Value = "Primera"
$ws.Range("J213").Value = 1700
$ws.Range("K213").Value = 14000
$ws.Range("L213").Value = 15000
$ws.Range("M213").Value = 14441
$ws.Range("N213").Value = '$/saco 25 kilos'
$ws.Range("O213").Value = "Región del Maule"
$ws.Range("P213").Value = 578
$ws.Range("Q213").Value = 25
$ws.Range("R213").Value = "Hortaliza"
